# SectorGroup.xlsx update:
# The "category" and "group" name/code columns get swapped so that
# columns D/E (names) exchange places, and columns F/G (codes)
# exchange places, for every row (including the header):
#   D: codeforiati:category-name  -> codeforiati:group-name
#   E: codeforiati:group-name     -> codeforiati:category-name
#   F: codeforiati:group-code     -> codeforiati:category-code
#   G: codeforiati:category-code  -> codeforiati:group-code

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value()
    $eVal = $ws.Cells.Item($r, 5).Value()
    $fVal = $ws.Cells.Item($r, 6).Value()
    $gVal = $ws.Cells.Item($r, 7).Value()

    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
    $ws.Cells.Item($r, 6).Value = $gVal
    $ws.Cells.Item($r, 7).Value = $fVal
}
